$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text looks like a plain number need an explicit
# text format, otherwise Excel COM auto-converts the assigned string into a
# numeric value (losing trailing zeros / exact text form).

$ws.Range("D2").Value = '37.445.09'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '2.049.68'
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.61'
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.612'
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.21'
$ws.Range("E8").Value = '  -3.58%  '
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("D12").Value = '2.352.55'
$ws.Range("E12").Value = '  -1.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.52'
$ws.Range("E13").Value = '  -4.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.58'
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("E15").Value = '  -3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.24'
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").Value = '2.046.23'
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").Value = '37.329.30'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.06'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.79'
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("D21").Value = '0.0₃0845'
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.62'
$ws.Range("E22").Value = '  -1.87%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  -4.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.49'
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.31'
$ws.Range("E27").Value = '  -1.94%  '
$ws.Range("E28").Value = '  -4.19%  '
$ws.Range("E29").Value = '  -2.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.90'
$ws.Range("E30").Value = '  -3.11%  '
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.53'
$ws.Range("E32").Value = '  -3.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0611'
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.53'
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("E35").Value = '  -3.74%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  -3.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  +0.76%  '
$ws.Range("E40").Value = '  -5.78%  '
$ws.Range("D41").Value = '1.499.69'
$ws.Range("E41").Value = '  +3.29%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.87'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.79'
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '95.89'
$ws.Range("E44").Value = '  -5.19%  '
$ws.Range("E45").Value = '  -3.63%  '
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("E47").Value = '  -4.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.23'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.93'
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.76'
$ws.Range("E50").Value = '  -8.51%  '
$ws.Range("D51").Value = '2.239.14'
$ws.Range("E51").Value = '  -1.80%  '
